$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(2,2).Value = "Wnt3"
$ws.Cells.Item(2,3).Value = "Fzd8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.2071343333333333
$ws.Cells.Item(2,8).Value = 0.6214029999999999
$ws.Cells.Item(2,9).Value = 0.4218324028717592
$ws.Cells.Item(2,10).Value = 0.4218324028717592
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.157506
$ws.Cells.Item(2,14).Value = 6.472517999999999
$ws.Cells.Item(2,15).Value = 0.3549648016839517
$ws.Cells.Item(2,16).Value = 0.3549648016839516
$ws.Cells.Item(2,17).Value = 0.4468935669726666
$ws.Cells.Item(2,18).Value = 4.022042102753999
$ws.Cells.Item(2,19).Value = 0.1497356552292388
$ws.Cells.Item(2,20).Value = 0.1497356552292388

# Row 3
$ws.Cells.Item(3,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(3,2).Value = "Wnt3"
$ws.Cells.Item(3,3).Value = "Fzd8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.2071343333333333
$ws.Cells.Item(3,8).Value = 0.6214029999999999
$ws.Cells.Item(3,9).Value = 0.4218324028717592
$ws.Cells.Item(3,10).Value = 0.4218324028717592
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.913654666666667
$ws.Cells.Item(3,14).Value = 8.740964
$ws.Cells.Item(3,15).Value = 0.4793705560628122
$ws.Cells.Item(3,16).Value = 0.4793705560628121
$ws.Cells.Item(3,17).Value = 0.6035179169435555
$ws.Cells.Item(3,18).Value = 5.431661252491999
$ws.Cells.Item(3,19).Value = 0.2022140335299474
$ws.Cells.Item(3,20).Value = 0.2022140335299474

# Row 4
$ws.Cells.Item(4,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,2).Value = "Wnt3"
$ws.Cells.Item(4,3).Value = "Fzd8"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.2071343333333333
$ws.Cells.Item(4,8).Value = 0.6214029999999999
$ws.Cells.Item(4,9).Value = 0.4218324028717592
$ws.Cells.Item(4,10).Value = 0.4218324028717592
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.018986
$ws.Cells.Item(4,14).Value = 0.05695799999999999
$ws.Cells.Item(4,15).Value = 0.003123681567871193
$ws.Cells.Item(4,16).Value = 0.003123681567871192
$ws.Cells.Item(4,17).Value = 0.003932652452666666
$ws.Cells.Item(4,18).Value = 0.03539387207399999
$ws.Cells.Item(4,19).Value = 0.00131767010158133
$ws.Cells.Item(4,20).Value = 0.001317670101581329

# Row 5
$ws.Cells.Item(5,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(5,2).Value = "Wnt3"
$ws.Cells.Item(5,3).Value = "Fzd8"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.2071343333333333
$ws.Cells.Item(5,8).Value = 0.6214029999999999
$ws.Cells.Item(5,9).Value = 0.4218324028717592
$ws.Cells.Item(5,10).Value = 0.4218324028717592
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.9848966666666668
$ws.Cells.Item(5,14).Value = 2.95469
$ws.Cells.Item(5,15).Value = 0.1620406385718132
$ws.Cells.Item(5,16).Value = 0.1620406385718132
$ws.Cells.Item(5,17).Value = 0.2040059144522222
$ws.Cells.Item(5,18).Value = 1.83605323007
$ws.Cells.Item(5,19).Value = 0.06835399193162224
$ws.Cells.Item(5,20).Value = 0.06835399193162223

# Row 6
$ws.Cells.Item(6,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value = "Wnt3"
$ws.Cells.Item(6,3).Value = "Fzd8"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.2071343333333333
$ws.Cells.Item(6,8).Value = 0.6214029999999999
$ws.Cells.Item(6,9).Value = 0.4218324028717592
$ws.Cells.Item(6,10).Value = 0.4218324028717592
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.003041
$ws.Cells.Item(6,14).Value = 0.009123
$ws.Cells.Item(6,15).Value = 0.0005003221135518961
$ws.Cells.Item(6,16).Value = 0.000500322113551896
$ws.Cells.Item(6,17).Value = 0.0006298955076666665
$ws.Cells.Item(6,18).Value = 0.005669059568999999
$ws.Cells.Item(6,19).Value = 0.0002110520793694735
$ws.Cells.Item(6,20).Value = 0.0002110520793694734

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Wnt3"
$ws.Cells.Item(7,3).Value = "Fzd8"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.01689933333333333
$ws.Cells.Item(7,8).Value = 0.050698
$ws.Cells.Item(7,9).Value = 0.0344157642637587
$ws.Cells.Item(7,10).Value = 0.0344157642637587
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.157506
$ws.Cells.Item(7,14).Value = 6.472517999999999
$ws.Cells.Item(7,15).Value = 0.3549648016839517
$ws.Cells.Item(7,16).Value = 0.3549648016839516
$ws.Cells.Item(7,17).Value = 0.03646041306266666
$ws.Cells.Item(7,18).Value = 0.328143717564
$ws.Cells.Item(7,19).Value = 0.01221638493668674
$ws.Cells.Item(7,20).Value = 0.01221638493668674

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Wnt3"
$ws.Cells.Item(8,3).Value = "Fzd8"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.01689933333333333
$ws.Cells.Item(8,8).Value = 0.050698
$ws.Cells.Item(8,9).Value = 0.0344157642637587
$ws.Cells.Item(8,10).Value = 0.0344157642637587
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.913654666666667
$ws.Cells.Item(8,14).Value = 8.740964
$ws.Cells.Item(8,15).Value = 0.4793705560628122
$ws.Cells.Item(8,16).Value = 0.4793705560628121
$ws.Cells.Item(8,17).Value = 0.04923882143022221
$ws.Cells.Item(8,18).Value = 0.443149392872
$ws.Cells.Item(8,19).Value = 0.01649790405244467
$ws.Cells.Item(8,20).Value = 0.01649790405244467

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Wnt3"
$ws.Cells.Item(9,3).Value = "Fzd8"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.01689933333333333
$ws.Cells.Item(9,8).Value = 0.050698
$ws.Cells.Item(9,9).Value = 0.0344157642637587
$ws.Cells.Item(9,10).Value = 0.0344157642637587
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.018986
$ws.Cells.Item(9,14).Value = 0.05695799999999999
$ws.Cells.Item(9,15).Value = 0.003123681567871193
$ws.Cells.Item(9,16).Value = 0.003123681567871192
$ws.Cells.Item(9,17).Value = 0.0003208507426666666
$ws.Cells.Item(9,18).Value = 0.002887656684
$ws.Cells.Item(9,19).Value = 0.0001075038884749032
$ws.Cells.Item(9,20).Value = 0.0001075038884749032

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Wnt3"
$ws.Cells.Item(10,3).Value = "Fzd8"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.01689933333333333
$ws.Cells.Item(10,8).Value = 0.050698
$ws.Cells.Item(10,9).Value = 0.0344157642637587
$ws.Cells.Item(10,10).Value = 0.0344157642637587
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.9848966666666668
$ws.Cells.Item(10,14).Value = 2.95469
$ws.Cells.Item(10,15).Value = 0.1620406385718132
$ws.Cells.Item(10,16).Value = 0.1620406385718132
$ws.Cells.Item(10,17).Value = 0.01664409706888889
$ws.Cells.Item(10,18).Value = 0.14979687362
$ws.Cells.Item(10,19).Value = 0.00557675241823645
$ws.Cells.Item(10,20).Value = 0.005576752418236449

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Wnt3"
$ws.Cells.Item(11,3).Value = "Fzd8"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.01689933333333333
$ws.Cells.Item(11,8).Value = 0.050698
$ws.Cells.Item(11,9).Value = 0.0344157642637587
$ws.Cells.Item(11,10).Value = 0.0344157642637587
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.003041
$ws.Cells.Item(11,14).Value = 0.009123
$ws.Cells.Item(11,15).Value = 0.0005003221135518961
$ws.Cells.Item(11,16).Value = 0.000500322113551896
$ws.Cells.Item(11,17).Value = 0.00005139087266666666
$ws.Cells.Item(11,18).Value = 0.000462517854
$ws.Cells.Item(11,19).Value = 0.00001721896791594757
$ws.Cells.Item(11,20).Value = 0.00001721896791594756

# Row 12
$ws.Cells.Item(12,1).Value = "Neutrophils"
$ws.Cells.Item(12,2).Value = "Wnt3"
$ws.Cells.Item(12,3).Value = "Fzd8"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.1481273333333333
$ws.Cells.Item(12,8).Value = 0.444382
$ws.Cells.Item(12,9).Value = 0.3016636978787648
$ws.Cells.Item(12,10).Value = 0.3016636978787648
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.157506
$ws.Cells.Item(12,14).Value = 6.472517999999999
$ws.Cells.Item(12,15).Value = 0.3549648016839517
$ws.Cells.Item(12,16).Value = 0.3549648016839516
$ws.Cells.Item(12,17).Value = 0.3195856104306666
$ws.Cells.Item(12,18).Value = 2.876270493876
$ws.Cells.Item(12,19).Value = 0.1070799946927833
$ws.Cells.Item(12,20).Value = 0.1070799946927833

# Row 13
$ws.Cells.Item(13,1).Value = "Neutrophils"
$ws.Cells.Item(13,2).Value = "Wnt3"
$ws.Cells.Item(13,3).Value = "Fzd8"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.1481273333333333
$ws.Cells.Item(13,8).Value = 0.444382
$ws.Cells.Item(13,9).Value = 0.3016636978787648
$ws.Cells.Item(13,10).Value = 0.3016636978787648
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 2.913654666666667
$ws.Cells.Item(13,14).Value = 8.740964
$ws.Cells.Item(13,15).Value = 0.4793705560628122
$ws.Cells.Item(13,16).Value = 0.4793705560628121
$ws.Cells.Item(13,17).Value = 0.4315918960275555
$ws.Cells.Item(13,18).Value = 3.884327064248
$ws.Cells.Item(13,19).Value = 0.1446086945961077
$ws.Cells.Item(13,20).Value = 0.1446086945961077

# Row 14
$ws.Cells.Item(14,1).Value = "Neutrophils"
$ws.Cells.Item(14,2).Value = "Wnt3"
$ws.Cells.Item(14,3).Value = "Fzd8"
$ws.Cells.Item(14,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.1481273333333333
$ws.Cells.Item(14,8).Value = 0.444382
$ws.Cells.Item(14,9).Value = 0.3016636978787648
$ws.Cells.Item(14,10).Value = 0.3016636978787648
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.018986
$ws.Cells.Item(14,14).Value = 0.05695799999999999
$ws.Cells.Item(14,15).Value = 0.003123681567871193
$ws.Cells.Item(14,16).Value = 0.003123681567871192
$ws.Cells.Item(14,17).Value = 0.002812345550666667
$ws.Cells.Item(14,18).Value = 0.025311109956
$ws.Cells.Item(14,19).Value = 0.0009423013327597621
$ws.Cells.Item(14,20).Value = 0.0009423013327597619

# Row 15
$ws.Cells.Item(15,1).Value = "Neutrophils"
$ws.Cells.Item(15,2).Value = "Wnt3"
$ws.Cells.Item(15,3).Value = "Fzd8"
$ws.Cells.Item(15,4).Value = "MuSCs"
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.1481273333333333
$ws.Cells.Item(15,8).Value = 0.444382
$ws.Cells.Item(15,9).Value = 0.3016636978787648
$ws.Cells.Item(15,10).Value = 0.3016636978787648
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.9848966666666668
$ws.Cells.Item(15,14).Value = 2.95469
$ws.Cells.Item(15,15).Value = 0.1620406385718132
$ws.Cells.Item(15,16).Value = 0.1620406385718132
$ws.Cells.Item(15,17).Value = 0.1458901168422222
$ws.Cells.Item(15,18).Value = 1.31301105158
$ws.Cells.Item(15,19).Value = 0.04888177823820959
$ws.Cells.Item(15,20).Value = 0.04888177823820958

# Row 16
$ws.Cells.Item(16,1).Value = "Neutrophils"
$ws.Cells.Item(16,2).Value = "Wnt3"
$ws.Cells.Item(16,3).Value = "Fzd8"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.1481273333333333
$ws.Cells.Item(16,8).Value = 0.444382
$ws.Cells.Item(16,9).Value = 0.3016636978787648
$ws.Cells.Item(16,10).Value = 0.3016636978787648
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.003041
$ws.Cells.Item(16,14).Value = 0.009123
$ws.Cells.Item(16,15).Value = 0.0005003221135518961
$ws.Cells.Item(16,16).Value = 0.000500322113551896
$ws.Cells.Item(16,17).Value = 0.0004504552206666666
$ws.Cells.Item(16,18).Value = 0.004054096986
$ws.Cells.Item(16,19).Value = 0.0001509290189045843
$ws.Cells.Item(16,20).Value = 0.0001509290189045842

# Row 17
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Wnt3"
$ws.Cells.Item(17,3).Value = "Fzd8"
$ws.Cells.Item(17,4).Value = "ECs"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 0.6666666666666666
$ws.Cells.Item(17,7).Value = 0.1188736666666667
$ws.Cells.Item(17,8).Value = 0.356621
$ws.Cells.Item(17,9).Value = 0.2420881349857173
$ws.Cells.Item(17,10).Value = 0.2420881349857172
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 2.157506
$ws.Cells.Item(17,14).Value = 6.472517999999999
$ws.Cells.Item(17,15).Value = 0.3549648016839517
$ws.Cells.Item(17,16).Value = 0.3549648016839516
$ws.Cells.Item(17,17).Value = 0.2564706490753333
$ws.Cells.Item(17,18).Value = 2.308235841678
$ws.Cells.Item(17,19).Value = 0.08593276682524285
$ws.Cells.Item(17,20).Value = 0.08593276682524283

# Row 18
$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Wnt3"
$ws.Cells.Item(18,3).Value = "Fzd8"
$ws.Cells.Item(18,4).Value = "FAPs"
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(18,6).Value = 0.6666666666666666
$ws.Cells.Item(18,7).Value = 0.1188736666666667
$ws.Cells.Item(18,8).Value = 0.356621
$ws.Cells.Item(18,9).Value = 0.2420881349857173
$ws.Cells.Item(18,10).Value = 0.2420881349857172
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 2.913654666666667
$ws.Cells.Item(18,14).Value = 8.740964
$ws.Cells.Item(18,15).Value = 0.4793705560628122
$ws.Cells.Item(18,16).Value = 0.4793705560628121
$ws.Cells.Item(18,17).Value = 0.3463568136271111
$ws.Cells.Item(18,18).Value = 3.117211322644
$ws.Cells.Item(18,19).Value = 0.1160499238843124
$ws.Cells.Item(18,20).Value = 0.1160499238843124

# Row 19
$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Wnt3"
$ws.Cells.Item(19,3).Value = "Fzd8"
$ws.Cells.Item(19,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19,5).Value = 2
$ws.Cells.Item(19,6).Value = 0.6666666666666666
$ws.Cells.Item(19,7).Value = 0.1188736666666667
$ws.Cells.Item(19,8).Value = 0.356621
$ws.Cells.Item(19,9).Value = 0.2420881349857173
$ws.Cells.Item(19,10).Value = 0.2420881349857172
$ws.Cells.Item(19,11).Value = 2
$ws.Cells.Item(19,12).Value = 0.6666666666666666
$ws.Cells.Item(19,13).Value = 0.018986
$ws.Cells.Item(19,14).Value = 0.05695799999999999
$ws.Cells.Item(19,15).Value = 0.003123681567871193
$ws.Cells.Item(19,16).Value = 0.003123681567871192
$ws.Cells.Item(19,17).Value = 0.002256935435333333
$ws.Cells.Item(19,18).Value = 0.020312418918
$ws.Cells.Item(19,19).Value = 0.0007562062450551984
$ws.Cells.Item(19,20).Value = 0.0007562062450551981

# Row 20
$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Wnt3"
$ws.Cells.Item(20,3).Value = "Fzd8"
$ws.Cells.Item(20,4).Value = "MuSCs"
$ws.Cells.Item(20,5).Value = 2
$ws.Cells.Item(20,6).Value = 0.6666666666666666
$ws.Cells.Item(20,7).Value = 0.1188736666666667
$ws.Cells.Item(20,8).Value = 0.356621
$ws.Cells.Item(20,9).Value = 0.2420881349857173
$ws.Cells.Item(20,10).Value = 0.2420881349857172
$ws.Cells.Item(20,11).Value = 3
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 0.9848966666666668
$ws.Cells.Item(20,14).Value = 2.95469
$ws.Cells.Item(20,15).Value = 0.1620406385718132
$ws.Cells.Item(20,16).Value = 0.1620406385718132
$ws.Cells.Item(20,17).Value = 0.1170782780544445
$ws.Cells.Item(20,18).Value = 1.05370450249
$ws.Cells.Item(20,19).Value = 0.03922811598374494
$ws.Cells.Item(20,20).Value = 0.03922811598374493

# Row 21
$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Wnt3"
$ws.Cells.Item(21,3).Value = "Fzd8"
$ws.Cells.Item(21,4).Value = "Resolving-Mac"
$ws.Cells.Item(21,5).Value = 2
$ws.Cells.Item(21,6).Value = 0.6666666666666666
$ws.Cells.Item(21,7).Value = 0.1188736666666667
$ws.Cells.Item(21,8).Value = 0.356621
$ws.Cells.Item(21,9).Value = 0.2420881349857173
$ws.Cells.Item(21,10).Value = 0.2420881349857172
$ws.Cells.Item(21,11).Value = 1
$ws.Cells.Item(21,12).Value = 0.3333333333333333
$ws.Cells.Item(21,13).Value = 0.003041
$ws.Cells.Item(21,14).Value = 0.009123
$ws.Cells.Item(21,15).Value = 0.0005003221135518961
$ws.Cells.Item(21,16).Value = 0.000500322113551896
$ws.Cells.Item(21,17).Value = 0.0003614948203333333
$ws.Cells.Item(21,18).Value = 0.003253453383
$ws.Cells.Item(21,19).Value = 0.0001211220473618908
$ws.Cells.Item(21,20).Value = 0.0001211220473618907
